$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 30.93275633333333
$ws.Range("H2").Value = 92.798269
$ws.Range("I2").Value = 0.6015903584115317
$ws.Range("J2").Value = 0.6015903584115317
$ws.Range("M2").Value = 0.03759866666666666
$ws.Range("O2").Value = 0.0008308188570719403
$ws.Range("P2").Value = 0.0008308188570719402
$ws.Range("Q2").Value = 1.163030394458222
$ws.Range("R2").Value = 10.467273550124
$ws.Range("S2").Value = 0.0004998126140009678
$ws.Range("T2").Value = 0.0004998126140009677
$ws.Range("G3").Value = 30.93275633333333
$ws.Range("H3").Value = 92.798269
$ws.Range("I3").Value = 0.6015903584115317
$ws.Range("J3").Value = 0.6015903584115317
$ws.Range("M3").Value = 21.57232766666667
$ws.Range("N3").Value = 64.716983
$ws.Range("O3").Value = 0.4766843669031189
$ws.Range("P3").Value = 0.4766843669031188
$ws.Range("Q3").Value = 667.2915552558252
$ws.Range("R3").Value = 6005.623997302428
$ws.Range("S3").Value = 0.2867687191344214
$ws.Range("T3").Value = 0.2867687191344214
$ws.Range("G4").Value = 30.93275633333333
$ws.Range("H4").Value = 92.798269
$ws.Range("I4").Value = 0.6015903584115317
$ws.Range("J4").Value = 0.6015903584115317
$ws.Range("O4").Value = 0.0007995515599313469
$ws.Range("P4").Value = 0.0007995515599313467
$ws.Range("Q4").Value = 1.119260544246556
$ws.Range("R4").Value = 10.073344898219
$ws.Range("S4").Value = 0.0004810025095075982
$ws.Range("T4").Value = 0.0004810025095075982
$ws.Range("G5").Value = 30.93275633333333
$ws.Range("H5").Value = 92.798269
$ws.Range("I5").Value = 0.6015903584115317
$ws.Range("J5").Value = 0.6015903584115317
$ws.Range("M5").Value = 23.608841
$ws.Range("N5").Value = 70.826523
$ws.Range("O5").Value = 0.5216852626798778
$ws.Range("P5").Value = 0.5216852626798777
$ws.Range("Q5").Value = 730.2865259654096
$ws.Range("R5").Value = 6572.578733688687
$ws.Range("S5").Value = 0.3138408241536018
$ws.Range("T5").Value = 0.3138408241536017
$ws.Range("I6").Value = 0.002688776579266707
$ws.Range("J6").Value = 0.002688776579266707
$ws.Range("M6").Value = 0.03759866666666666
$ws.Range("O6").Value = 0.0008308188570719403
$ws.Range("P6").Value = 0.0008308188570719402
$ws.Range("Q6").Value = 0.005198103396888889
$ws.Range("R6").Value = 0.046782930572
$ws.Range("S6").Value = 0.000002233886284508167
$ws.Range("T6").Value = 0.000002233886284508166
$ws.Range("I7").Value = 0.002688776579266707
$ws.Range("J7").Value = 0.002688776579266707
$ws.Range("M7").Value = 21.57232766666667
$ws.Range("N7").Value = 64.716983
$ws.Range("O7").Value = 0.4766843669031189
$ws.Range("P7").Value = 0.4766843669031188
$ws.Range("Q7").Value = 2.982424635347889
$ws.Range("R7").Value = 26.841821718131
$ws.Range("S7").Value = 0.001281697761431684
$ws.Range("T7").Value = 0.001281697761431684
$ws.Range("I8").Value = 0.002688776579266707
$ws.Range("J8").Value = 0.002688776579266707
$ws.Range("O8").Value = 0.0007995515599313469
$ws.Range("P8").Value = 0.0007995515599313467
$ws.Range("R8").Value = 0.04502228710699999
$ws.Range("S8").Value = 0.000002149815508259566
$ws.Range("T8").Value = 0.000002149815508259566
$ws.Range("I9").Value = 0.002688776579266707
$ws.Range("J9").Value = 0.002688776579266707
$ws.Range("M9").Value = 23.608841
$ws.Range("N9").Value = 70.826523
$ws.Range("O9").Value = 0.5216852626798778
$ws.Range("P9").Value = 0.5216852626798777
$ws.Range("Q9").Value = 3.263977355545666
$ws.Range("R9").Value = 29.375796199911
$ws.Range("S9").Value = 0.001402695116042255
$ws.Range("T9").Value = 0.001402695116042255
$ws.Range("G10").Value = 13.65672433333333
$ws.Range("H10").Value = 40.970173
$ws.Range("I10").Value = 0.2656004398018724
$ws.Range("J10").Value = 0.2656004398018724
$ws.Range("M10").Value = 0.03759866666666666
$ws.Range("O10").Value = 0.0008308188570719403
$ws.Range("P10").Value = 0.0008308188570719402
$ws.Range("Q10").Value = 0.5134746259675556
$ws.Range("R10").Value = 4.621271633708
$ws.Range("S10").Value = 0.0002206658538339963
$ws.Range("T10").Value = 0.0002206658538339963
$ws.Range("G11").Value = 13.65672433333333
$ws.Range("H11").Value = 40.970173
$ws.Range("I11").Value = 0.2656004398018724
$ws.Range("J11").Value = 0.2656004398018724
$ws.Range("M11").Value = 21.57232766666667
$ws.Range("N11").Value = 64.716983
$ws.Range("O11").Value = 0.4766843669031189
$ws.Range("P11").Value = 0.4766843669031188
$ws.Range("Q11").Value = 294.6073321720066
$ws.Range("R11").Value = 2651.465989548059
$ws.Range("S11").Value = 0.1266075774961455
$ws.Range("T11").Value = 0.1266075774961455
$ws.Range("G12").Value = 13.65672433333333
$ws.Range("H12").Value = 40.970173
$ws.Range("I12").Value = 0.2656004398018724
$ws.Range("J12").Value = 0.2656004398018724
$ws.Range("O12").Value = 0.0007995515599313469
$ws.Range("P12").Value = 0.0007995515599313467
$ws.Range("Q12").Value = 0.4941503610358889
$ws.Range("R12").Value = 4.447353249323
$ws.Range("S12").Value = 0.0002123612459620389
$ws.Range("T12").Value = 0.0002123612459620388
$ws.Range("G13").Value = 13.65672433333333
$ws.Range("H13").Value = 40.970173
$ws.Range("I13").Value = 0.2656004398018724
$ws.Range("J13").Value = 0.2656004398018724
$ws.Range("M13").Value = 23.608841
$ws.Range("N13").Value = 70.826523
$ws.Range("O13").Value = 0.5216852626798778
$ws.Range("P13").Value = 0.5216852626798777
$ws.Range("Q13").Value = 322.4194333664977
$ws.Range("R13").Value = 2901.774900298479
$ws.Range("S13").Value = 0.1385598352059309
$ws.Range("T13").Value = 0.1385598352059308
$ws.Range("G14").Value = 0.1180373333333333
$ws.Range("H14").Value = 0.354112
$ws.Range("I14").Value = 0.002295628650118725
$ws.Range("J14").Value = 0.002295628650118725
$ws.Range("M14").Value = 0.03759866666666666
$ws.Range("O14").Value = 0.0008308188570719403
$ws.Range("P14").Value = 0.0008308188570719402
$ws.Range("Q14").Value = 0.004438046350222221
$ws.Range("R14").Value = 0.039942417152
$ws.Range("S14").Value = 0.00000190725157135324
$ws.Range("T14").Value = 0.00000190725157135324
$ws.Range("G15").Value = 0.1180373333333333
$ws.Range("H15").Value = 0.354112
$ws.Range("I15").Value = 0.002295628650118725
$ws.Range("J15").Value = 0.002295628650118725
$ws.Range("M15").Value = 21.57232766666667
$ws.Range("N15").Value = 64.716983
$ws.Range("O15").Value = 0.4766843669031189
$ws.Range("P15").Value = 0.4766843669031188
$ws.Range("Q15").Value = 2.546340031566222
$ws.Range("R15").Value = 22.917060284096
$ws.Range("S15").Value = 0.001094290289726506
$ws.Range("T15").Value = 0.001094290289726506
$ws.Range("G16").Value = 0.1180373333333333
$ws.Range("H16").Value = 0.354112
$ws.Range("I16").Value = 0.002295628650118725
$ws.Range("J16").Value = 0.002295628650118725
$ws.Range("O16").Value = 0.0007995515599313469
$ws.Range("P16").Value = 0.0007995515599313467
$ws.Range("Q16").Value = 0.004271023523555554
$ws.Range("R16").Value = 0.03843921171199999
$ws.Range("S16").Value = 0.000001835473468225519
$ws.Range("T16").Value = 0.000001835473468225518
$ws.Range("G17").Value = 0.1180373333333333
$ws.Range("H17").Value = 0.354112
$ws.Range("I17").Value = 0.002295628650118725
$ws.Range("J17").Value = 0.002295628650118725
$ws.Range("M17").Value = 23.608841
$ws.Range("N17").Value = 70.826523
$ws.Range("O17").Value = 0.5216852626798778
$ws.Range("P17").Value = 0.5216852626798777
$ws.Range("Q17").Value = 2.786724634730666
$ws.Range("R17").Value = 25.080521712576
$ws.Range("S17").Value = 0.00119759563535264
$ws.Range("T17").Value = 0.00119759563535264
$ws.Range("G18").Value = 6.572534333333333
$ws.Range("H18").Value = 19.717603
$ws.Range("I18").Value = 0.1278247965572105
$ws.Range("J18").Value = 0.1278247965572105
$ws.Range("M18").Value = 0.03759866666666666
$ws.Range("O18").Value = 0.0008308188570719403
$ws.Range("P18").Value = 0.0008308188570719402
$ws.Range("Q18").Value = 0.2471185275542222
$ws.Range("R18").Value = 2.224066747988
$ws.Range("S18").Value = 0.0001061992513811149
$ws.Range("T18").Value = 0.0001061992513811149
$ws.Range("G19").Value = 6.572534333333333
$ws.Range("H19").Value = 19.717603
$ws.Range("I19").Value = 0.1278247965572105
$ws.Range("J19").Value = 0.1278247965572105
$ws.Range("M19").Value = 21.57232766666667
$ws.Range("N19").Value = 64.716983
$ws.Range("O19").Value = 0.4766843669031189
$ws.Range("P19").Value = 0.4766843669031188
$ws.Range("Q19").Value = 141.7848642390832
$ws.Range("R19").Value = 1276.063778151749
$ws.Range("S19").Value = 0.06093208222139385
$ws.Range("T19").Value = 0.06093208222139385
$ws.Range("G20").Value = 6.572534333333333
$ws.Range("H20").Value = 19.717603
$ws.Range("I20").Value = 0.1278247965572105
$ws.Range("J20").Value = 0.1278247965572105
$ws.Range("O20").Value = 0.0007995515599313469
$ws.Range("P20").Value = 0.0007995515599313467
$ws.Range("Q20").Value = 0.2378183914725555
$ws.Range("R20").Value = 2.140365523253
$ws.Range("S20").Value = 0.0001022025154852247
$ws.Range("T20").Value = 0.0001022025154852247
$ws.Range("G21").Value = 6.572534333333333
$ws.Range("H21").Value = 19.717603
$ws.Range("I21").Value = 0.1278247965572105
$ws.Range("J21").Value = 0.1278247965572105
$ws.Range("M21").Value = 23.608841
$ws.Range("N21").Value = 70.826523
$ws.Range("O21").Value = 0.5216852626798778
$ws.Range("P21").Value = 0.5216852626798777
$ws.Range("Q21").Value = 155.1699180427076
$ws.Range("R21").Value = 1396.529262384369
$ws.Range("S21").Value = 0.0666843125689503
$ws.Range("T21").Value = 0.06668431256895029
